# Resolve the empty/duplicate footnote bug:
# - footnote id=116 is a spurious empty footnote (just "།") with no real
#   content; both its body reference and its definition must be removed.
# - footnote id=22's body text ("།" only) was actually meant to contain
#   a real note; restore the intended text.
# - footnote id=30's body text had a stray trailing "a" typo; drop it.

$d = $word.ActiveDocument

$emptyText = "།"

# Find + delete the trailing empty footnote (the last footnote whose entire
# body text -- apart from the separating space run -- is just the
# lone punctuation mark "།").  Walking backwards guarantees we hit the
# *last* such footnote, which is the spurious one introduced by the bug
# (there is an earlier, legitimate one -- id=22 -- with the same
# placeholder text that must be repaired instead of removed).
$targetIndex = -1
for ($i = $d.Footnotes.Count; $i -ge 1; $i--) {
    $fn = $d.Footnotes.Item($i)
    if ($fn.Range.Text.Trim() -eq $emptyText) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $d.Footnotes.Item($targetIndex).Delete()
}

# Repair footnote id=22 (now still the 2nd footnote): give it its real text.
$fn22 = $d.Footnotes.Item(2)
$fn22.Range.Text = "འཕགས་པ། ཞེས་པར་མ་གཞན་ནང་མེད།"

# Footnote id=30 (10th footnote): drop the stray trailing "a".
$fn30 = $d.Footnotes.Item(10)
$fn30.Range.Text = "དེ་ནས། པེ་ཅིན།"
